$d = $word.ActiveDocument

# 1. "4.12.2018" -> "4.12"
$null = $d.Content.Find.Execute("4.12.2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4.12", 2)

# 2. Remove the _GoBack bookmark from its current location (end of the
#    "...made git repository." paragraph) - it will be re-added after the
#    new last paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Insert the four new paragraphs after the "made git repository." paragraph.
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.Text = "5.12"

$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "Created cnr.html and cnr.css files, started working on the number base conversion program"

$d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.Text = "6.12"

$d.Paragraphs(7).Range.InsertParagraphAfter()
# Trailing sentinel "X" lets us park the collapsed bookmark one char before
# the paragraph mark without hitting the "position == paragraph.End-1"
# placement bug, then we delete the sentinel so the bookmark ends up
# exactly where it needs to be (right after the real text, before the
# paragraph mark) without the run getting split around it.
$d.Paragraphs(8).Range.Text = "Got number base conversion to working order, worked more on the conversion page UI. X"

$p8 = $d.Paragraphs(8)
$bmPos = $p8.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelPos = $d.Paragraphs(8).Range.End - 1
$sentinelRange = $d.Range($sentinelPos - 1, $sentinelPos)
$sentinelRange.Delete()
